# "add test random packing spheres"
#
# Sheet1 ("Sheet1"): row 1 (B1:J1) becomes a text header row (shared
# strings) instead of the small integer codes, and column H / column J
# get swapped for the data rows (for rows 3-4 only column H had a value,
# so it simply moves over to column J). Sheet2 ("Sheet2") gets the same
# H/J swap applied to its row 1 and data rows.

$wb = $excel.ActiveWorkbook

# Cosmetic: the tab-bar/horizontal-scrollbar splitter ratio moved
# slightly (986 -> 989 in the raw XML, i.e. 0.986 -> 0.989).
$wb.Windows.Item(1).TabRatio = 0.989

# ---------------------------------------------------------------------
# Sheet1
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()

# New text headers (row 1, B:J) -- written left-to-right so the shared
# string table indices come out as 0..8 in this order.
$ws1.Range("B1").Value2 = "SBB-MRT"
$ws1.Range("C1").Value2 = "LIBB-MRT"
$ws1.Range("D1").Value2 = "QIBB-MRT"
$ws1.Range("E1").Value2 = "MR-MRT"
$ws1.Range("F1").Value2 = "CLI-MRT"
$ws1.Range("G1").Value2 = "PSM-MRT-A"
$ws1.Range("H1").Value2 = "PSM-MRT-B"
$ws1.Range("I1").Value2 = "IBM-MRT-A"
$ws1.Range("J1").Value2 = "IBM-MRT-B"

# Rows 3-4: only column H was populated: the value moves over to column
# J and the old H cell becomes empty.
for ($r = 3; $r -le 4; $r++) {
    $hCell = $ws1.Range("H$r")
    $hVal = $hCell.Value2
    $ws1.Range("J$r").Value2 = $hVal
    $hCell.Clear()
}

# Rows 5-12: both H and J are populated -- swap the two values.
for ($r = 5; $r -le 12; $r++) {
    $hCell = $ws1.Range("H$r")
    $jCell = $ws1.Range("J$r")
    $hVal = $hCell.Value2
    $jVal = $jCell.Value2
    $hCell.Value2 = $jVal
    $jCell.Value2 = $hVal
}

# Cosmetic: default column width / selection (best effort).
$ws1.StandardWidth = 8.23469387755102
$ws1.Range("J1:J12").Select()
$ws1.Range("K14").Activate()

# ---------------------------------------------------------------------
# Sheet2
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Activate()

# Row 1: swap H1 / J1.
$h1 = $ws2.Range("H1")
$j1 = $ws2.Range("J1")
$h1Val = $h1.Value2
$j1Val = $j1.Value2
$h1.Value2 = $j1Val
$j1.Value2 = $h1Val

# Rows 5-12: swap H / J.
for ($r = 5; $r -le 12; $r++) {
    $hCell = $ws2.Range("H$r")
    $jCell = $ws2.Range("J$r")
    $hVal = $hCell.Value2
    $jVal = $jCell.Value2
    $hCell.Value2 = $jVal
    $jCell.Value2 = $hVal
}

$ws2.Range("J1:J12").Select()

# Re-activate Sheet2 so it remains the workbook's active tab (matches
# the original file, where Sheet2 is tabSelected/activeTab).
$ws2.Activate()
